{"js": "// \"cambio en el paso 15, elimino la -u\"\n//\n// 1) Before the step that says \"para subir el archivo a Github hago git\n//    remote add origin ...\" insert three new numbered-list steps:\n//      \"git status \", \"git add .\", \"git commit -m (1 commit)\"\n// 2) Merge the three runs that make up \"para subir el archivo a Github \" +\n//    \"hago \" + \"git remote add origin \" into a single run (same text,\n//    concatenated) while leaving the following hyperlink run untouched.\n// 3) At the very end of the document append five new numbered-list steps:\n//      \"modifico el archivo de nuevo \", \"git status\", \"git add .\",\n//      \"git commit -m (3 commit)\", \"git push origin master \"\n//    (this is the \"elimino la -u\" part: the old final step used\n//    \"git push -u origin master\"; the new last step uses\n//    \"git push origin master \" without \"-u\").\n\nconst body = context.document.body;\n\n// --- Step 1: locate the \"para subir el archivo a Github\" step -----------\nconst target = body.search(\"para subir el archivo a Github\", { matchCase: false });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error('Could not find paragraph \"para subir el archivo a Github\"');\n}\n\nconst targetParagraphs = target.items[0].paragraphs;\ntargetParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = targetParagraphs.items[0];\n\n// Insert the three new steps immediately before it, in reading order.\ntargetParagraph.insertParagraph(\"git status \", Word.InsertLocation.before);\ntargetParagraph.insertParagraph(\"git add .\", Word.InsertLocation.before);\ntargetParagraph.insertParagraph(\"git commit -m (1 commit)\", Word.InsertLocation.before);\nawait context.sync();\n\n// --- Step 2: merge the 3 runs into a single run ---------------------------\nconst mergedText =\n  \"para subir el archivo a Github hago git remote add origin \";\nconst mergeRange = body.search(mergedText, { matchCase: false });\nmergeRange.load(\"items\");\nawait context.sync();\n\nif (mergeRange.items.length === 0) {\n  throw new Error(\"Could not find the text to merge\");\n}\nmergeRange.items[0].insertText(mergedText, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 3: append five new steps at the end of the document ------------\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nlet lastParagraph = allParagraphs.items[allParagraphs.items.length - 1];\n\nconst newSteps = [\n  \"modifico el archivo de nuevo \",\n  \"git status\",\n  \"git add .\",\n  \"git commit -m (3 commit)\",\n  \"git push origin master \",\n];\n\nfor (const stepText of newSteps) {\n  lastParagraph = lastParagraph.insertParagraph(stepText, Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# \"cambio en el paso 15, elimino la -u\"\n#\n# 1) Before the step that reads \"para subir el archivo a Github hago git\n#    remote add origin ...\" insert three new numbered-list steps:\n#      \"git status \", \"git add .\", \"git commit -m (1 commit)\"\n# 2) Merge the three runs that make up \"para subir el archivo a Github \" +\n#    \"hago \" + \"git remote add origin \" into a single run (identical\n#    concatenated text) while leaving the following hyperlink run untouched.\n# 3) At the very end of the document append five new numbered-list steps:\n#      \"modifico el archivo de nuevo \", \"git status\", \"git add .\",\n#      \"git commit -m (3 commit)\", \"git push origin master \"\n#    (the \"elimino la -u\" part: the old last step used\n#    \"git push -u origin master\"; the new final step drops \"-u\").\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the three new steps right before \"para subir...\" -----\n$newStepsBefore = @(\"git status \", \"git add .\", \"git commit -m (1 commit)\")\nforeach ($stepText in $newStepsBefore) {\n    $rng = $d.Content\n    $rng.Find.Execute(\"para subir el archivo a Github\") | Out-Null\n    $targetPara = $rng.Paragraphs(1)\n    $targetPara.Range.InsertParagraphBefore()\n    $targetPara.Range.Text = $stepText\n}\n\n# --- Step 2: merge the 3 runs into a single run ---------------------------\n$mergedText = \"para subir el archivo a Github hago git remote add origin \"\n$rng = $d.Content\n$rng.Find.Execute($mergedText, $false, $false, $false, $false, $false, $true, 1, $false, $mergedText, 1) | Out-Null\n\n# --- Step 3: append five new steps at the end of the document ------------\n$newStepsAfter = @(\n    \"modifico el archivo de nuevo \",\n    \"git status\",\n    \"git add .\",\n    \"git commit -m (3 commit)\",\n    \"git push origin master \"\n)\nforeach ($stepText in $newStepsAfter) {\n    $count = $d.Paragraphs.Count\n    $last = $d.Paragraphs($count)\n    $last.Range.InsertParagraphAfter()\n    $newLast = $d.Paragraphs($d.Paragraphs.Count)\n    $newLast.Range.Text = $stepText\n}\n"}
